$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 from "Good Morning" to "Good Morning11"
$ws.Range("E8").Value = "Good Morning11"
